$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (default, unformatted) cell used to restore style after
# forcing a text NumberFormat, so the written cells keep their original
# (unstyled) appearance instead of picking up a new @ text style index.
$defaultStyle = $ws.Range("B2").Style

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to Text format before assigning so that numeric-
    # looking strings (e.g. "1.00", "0.0000280") are kept verbatim as
    # text instead of being coerced into a floating point number.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $defaultStyle
}

Set-TextValue $ws.Range("D2") "65.547.13"
Set-TextValue $ws.Range("E2") "  +3.20%  "
Set-TextValue $ws.Range("D3") "3.407.49"
Set-TextValue $ws.Range("E3") "  +2.45%  "
Set-TextValue $ws.Range("E4") "  -0.04%  "
Set-TextValue $ws.Range("D5") "562.79"
Set-TextValue $ws.Range("E5") "  +3.17%  "
Set-TextValue $ws.Range("D6") "176.18"
Set-TextValue $ws.Range("E6") "  +2.72%  "
Set-TextValue $ws.Range("E7") "  +3.14%  "
Set-TextValue $ws.Range("D8") "3.397.04"
Set-TextValue $ws.Range("E8") "  +2.30%  "
Set-TextValue $ws.Range("D9") "0.999"
Set-TextValue $ws.Range("E9") "  -0.14%  "
Set-TextValue $ws.Range("D10") "0.171"
Set-TextValue $ws.Range("E10") "  +13.35%  "
Set-TextValue $ws.Range("E11") "  +3.81%  "
Set-TextValue $ws.Range("D12") "55.04"
Set-TextValue $ws.Range("E12") "  +3.30%  "
Set-TextValue $ws.Range("D13") "0.0000280"
Set-TextValue $ws.Range("E13") "  +6.27%  "
Set-TextValue $ws.Range("E14") "  +3.30%  "
Set-TextValue $ws.Range("D15") "3.950.36"
Set-TextValue $ws.Range("E15") "  +2.12%  "
Set-TextValue $ws.Range("D16") "18.38"
Set-TextValue $ws.Range("E16") "  +3.16%  "
Set-TextValue $ws.Range("D17") "3.405.42"
Set-TextValue $ws.Range("E17") "  +2.89%  "
Set-TextValue $ws.Range("E18") "  +2.20%  "
Set-TextValue $ws.Range("D19") "65.340.94"
Set-TextValue $ws.Range("E19") "  +2.87%  "
Set-TextValue $ws.Range("D20") "11.93"
Set-TextValue $ws.Range("E20") "  +2.25%  "
Set-TextValue $ws.Range("E21") "  +2.71%  "
Set-TextValue $ws.Range("D22") "471.00"
Set-TextValue $ws.Range("E22") "  +14.31%  "
Set-TextValue $ws.Range("D23") "5.11"
Set-TextValue $ws.Range("E23") "  +18.50%  "
Set-TextValue $ws.Range("D24") "4.16"
Set-TextValue $ws.Range("E24") "  +3.38%  "
Set-TextValue $ws.Range("D25") "86.74"
Set-TextValue $ws.Range("E25") "  +4.58%  "
Set-TextValue $ws.Range("D26") "13.43"
Set-TextValue $ws.Range("E26") "  -0.95%  "
Set-TextValue $ws.Range("D27") "10.93"
Set-TextValue $ws.Range("E27") "  +3.71%  "
Set-TextValue $ws.Range("E28") "  +7.06%  "
Set-TextValue $ws.Range("D29") "8.90"
Set-TextValue $ws.Range("E29") "  +4.28%  "
Set-TextValue $ws.Range("D30") "31.24"
Set-TextValue $ws.Range("E30") "  +7.86%  "
Set-TextValue $ws.Range("E31") "  +5.99%  "
Set-TextValue $ws.Range("E32") "  +2.65%  "
Set-TextValue $ws.Range("D33") "62.84"
Set-TextValue $ws.Range("E33") "  +9.15%  "
Set-TextValue $ws.Range("D34") "573.77"
Set-TextValue $ws.Range("E34") "  +0.05%  "
Set-TextValue $ws.Range("E35") "  +2.74%  "
Set-TextValue $ws.Range("E36") "  +0.01%  "
Set-TextValue $ws.Range("E37") "  -4.18%  "
Set-TextValue $ws.Range("D38") "3.54"
Set-TextValue $ws.Range("E38") "  +4.46%  "
Set-TextValue $ws.Range("D39") "35.97"
Set-TextValue $ws.Range("E39") "  +2.90%  "
Set-TextValue $ws.Range("E40") "  +3.70%  "
Set-TextValue $ws.Range("E41") "  +2.61%  "
Set-TextValue $ws.Range("D42") "3.092.00"
Set-TextValue $ws.Range("E42") "  -1.39%  "
Set-TextValue $ws.Range("D43") "1.00"
Set-TextValue $ws.Range("E43") "  -0.03%  "
Set-TextValue $ws.Range("E44") "  +2.58%  "
Set-TextValue $ws.Range("D45") "0.0418"
Set-TextValue $ws.Range("E45") "  +4.83%  "
Set-TextValue $ws.Range("D46") "2.51"
Set-TextValue $ws.Range("E46") "  +4.53%  "
Set-TextValue $ws.Range("D47") "0.135"
Set-TextValue $ws.Range("E47") "  +6.45%  "
Set-TextValue $ws.Range("D48") "3.16"
Set-TextValue $ws.Range("E48") "  -2.13%  "
Set-TextValue $ws.Range("E49") "  +0.32%  "
Set-TextValue $ws.Range("B50") "Monero"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D50") "136.68"
Set-TextValue $ws.Range("E50") "  +3.65%  "
Set-TextValue $ws.Range("B51") "THORChain"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D51") "8.36"
Set-TextValue $ws.Range("E51") "  +4.51%  "
